$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 0.2011661807580175
$ws.Cells.Item(2, 3).Value = 0.5276967930029155
$ws.Cells.Item(2, 10).Value = 0.02915451895043732
$ws.Cells.Item(2, 16).Value = 0.1545189504373178
$ws.Cells.Item(2, 19).Value = 0.08746355685131195
$ws.Cells.Item(3, 3).Value = 0.02659574468085106
$ws.Cells.Item(3, 10).Value = 0.05851063829787234
$ws.Cells.Item(3, 16).Value = 0.723404255319149
$ws.Cells.Item(3, 19).Value = 0.1914893617021277
$ws.Cells.Item(4, 10).Value = 0.07317073170731707
$ws.Cells.Item(4, 16).Value = 0.5853658536585366
$ws.Cells.Item(4, 19).Value = 0.3414634146341464
$ws.Cells.Item(6, 2).Value = 0.1231527093596059
$ws.Cells.Item(6, 4).Value = 0.009852216748768473
$ws.Cells.Item(6, 6).Value = 0.04926108374384237
$ws.Cells.Item(6, 10).Value = 0.2463054187192118
$ws.Cells.Item(6, 15).Value = 0.03448275862068965
$ws.Cells.Item(6, 17).Value = 0.1231527093596059
$ws.Cells.Item(6, 18).Value = 0.0541871921182266
$ws.Cells.Item(6, 19).Value = 0.3596059113300493
$ws.Cells.Item(7, 2).Value = 0.1382488479262673
$ws.Cells.Item(7, 4).Value = 0.02764976958525346
$ws.Cells.Item(7, 5).Value = 0.004608294930875576
$ws.Cells.Item(7, 6).Value = 0.05529953917050692
$ws.Cells.Item(7, 10).Value = 0.1658986175115207
$ws.Cells.Item(7, 15).Value = 0.02304147465437788
$ws.Cells.Item(7, 17).Value = 0.1474654377880184
$ws.Cells.Item(7, 18).Value = 0.05990783410138249
$ws.Cells.Item(7, 19).Value = 0.3778801843317972
$ws.Cells.Item(8, 2).Value = 0.08446866485013624
$ws.Cells.Item(8, 4).Value = 0.01362397820163488
$ws.Cells.Item(8, 6).Value = 0.06539509536784741
$ws.Cells.Item(8, 10).Value = 0.1198910081743869
$ws.Cells.Item(8, 15).Value = 0.0217983651226158
$ws.Cells.Item(8, 17).Value = 0.1416893732970027
$ws.Cells.Item(8, 18).Value = 0.08446866485013624
$ws.Cells.Item(8, 19).Value = 0.4686648501362398
$ws.Cells.Item(9, 2).Value = 0.07027027027027027
$ws.Cells.Item(9, 4).Value = 0.02702702702702703
$ws.Cells.Item(9, 6).Value = 0.03243243243243243
$ws.Cells.Item(9, 10).Value = 0.1837837837837838
$ws.Cells.Item(9, 17).Value = 0.1783783783783784
$ws.Cells.Item(9, 18).Value = 0.05405405405405406
$ws.Cells.Item(9, 19).Value = 0.4540540540540541
$ws.Cells.Item(10, 2).Value = 0.1328878990348923
$ws.Cells.Item(10, 4).Value = 0.01707498144023756
$ws.Cells.Item(10, 5).Value = 0.001484780994803266
$ws.Cells.Item(10, 6).Value = 0.066815144766147
$ws.Cells.Item(10, 10).Value = 0.1210096510764662
$ws.Cells.Item(10, 15).Value = 0.0178173719376392
$ws.Cells.Item(10, 17).Value = 0.2048997772828508
$ws.Cells.Item(10, 18).Value = 0.07720861172976985
$ws.Cells.Item(10, 19).Value = 0.3608017817371937
$ws.Cells.Item(11, 7).Value = 0.1481481481481481
$ws.Cells.Item(11, 10).Value = 0.09567901234567901
$ws.Cells.Item(11, 11).Value = 0.1944444444444444
$ws.Cells.Item(11, 12).Value = 0.5401234567901234
$ws.Cells.Item(11, 19).Value = 0.02160493827160494
$ws.Cells.Item(12, 7).Value = 0.7437185929648241
$ws.Cells.Item(12, 10).Value = 0.1658291457286432
$ws.Cells.Item(12, 11).Value = 0.02512562814070352
$ws.Cells.Item(12, 12).Value = 0.04020100502512563
$ws.Cells.Item(12, 19).Value = 0.02512562814070352
$ws.Cells.Item(13, 7).Value = 0.8484848484848485
$ws.Cells.Item(13, 10).Value = 0.1212121212121212
$ws.Cells.Item(13, 19).Value = 0.0303030303030303
$ws.Cells.Item(14, 7).Value = 0.8333333333333334
$ws.Cells.Item(14, 10).Value = 0.1666666666666667
$ws.Cells.Item(15, 6).Value = 0.04504504504504504
$ws.Cells.Item(15, 8).Value = 0.1486486486486487
$ws.Cells.Item(15, 9).Value = 0.07207207207207207
$ws.Cells.Item(15, 10).Value = 0.4099099099099099
$ws.Cells.Item(15, 11).Value = 0.03603603603603604
$ws.Cells.Item(15, 15).Value = 0.04504504504504504
$ws.Cells.Item(15, 19).Value = 0.2432432432432433
$ws.Cells.Item(16, 6).Value = 0.02010050251256281
$ws.Cells.Item(16, 8).Value = 0.1809045226130653
$ws.Cells.Item(16, 9).Value = 0.07537688442211055
$ws.Cells.Item(16, 10).Value = 0.4723618090452261
$ws.Cells.Item(16, 11).Value = 0.1105527638190955
$ws.Cells.Item(16, 15).Value = 0.02512562814070352
$ws.Cells.Item(16, 19).Value = 0.1155778894472362
$ws.Cells.Item(17, 6).Value = 0.007263922518159807
$ws.Cells.Item(17, 8).Value = 0.1162227602905569
$ws.Cells.Item(17, 9).Value = 0.1089588377723971
$ws.Cells.Item(17, 10).Value = 0.4769975786924939
$ws.Cells.Item(17, 11).Value = 0.09927360774818401
$ws.Cells.Item(17, 13).Value = 0.01452784503631961
$ws.Cells.Item(17, 15).Value = 0.05326876513317191
$ws.Cells.Item(17, 19).Value = 0.1234866828087167
$ws.Cells.Item(18, 6).Value = 0.005988023952095809
$ws.Cells.Item(18, 8).Value = 0.125748502994012
$ws.Cells.Item(18, 9).Value = 0.0658682634730539
$ws.Cells.Item(18, 10).Value = 0.4610778443113773
$ws.Cells.Item(18, 11).Value = 0.1137724550898204
$ws.Cells.Item(18, 13).Value = 0.01796407185628742
$ws.Cells.Item(18, 14).Value = 0.005988023952095809
$ws.Cells.Item(18, 15).Value = 0.08982035928143713
$ws.Cells.Item(18, 19).Value = 0.1137724550898204
$ws.Cells.Item(19, 6).Value = 0.01195219123505976
$ws.Cells.Item(19, 8).Value = 0.1816733067729084
$ws.Cells.Item(19, 9).Value = 0.07250996015936255
$ws.Cells.Item(19, 10).Value = 0.3920318725099601
$ws.Cells.Item(19, 11).Value = 0.1258964143426295
$ws.Cells.Item(19, 13).Value = 0.01832669322709163
$ws.Cells.Item(19, 14).Value = 0.00398406374501992
$ws.Cells.Item(19, 15).Value = 0.08207171314741035
$ws.Cells.Item(19, 19).Value = 0.1115537848605578

